$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 330.1111  # H18 was 373.14285
$ws.Cells.Item(18, 9).Value = 192.75  # I18 was 204
$ws.Cells.Item(18, 10).Value = 440  # J18 was 500
$ws.Cells.Item(18, 11).Value = 192.75  # K18 was 204
$ws.Cells.Item(18, 12).Value = 440  # L18 was 500
$ws.Cells.Item(18, 13).Value = 91.25  # M18 was 80
$ws.Cells.Item(18, 14).Value = -1008  # N18 was -1068
$ws.Cells.Item(64, 8).Value = 4009.1667  # H64 was 4211
$ws.Cells.Item(64, 10).Value = 4518.3335  # J64 was 5277.5
$ws.Cells.Item(64, 12).Value = 4518.3335  # L64 was 5277.5
$ws.Cells.Item(64, 14).Value = -5014.3335  # N64 was -5773.5
$ws.Cells.Item(67, 8).Value = 4009.1667  # H67 was 4211
$ws.Cells.Item(67, 10).Value = 4518.3335  # J67 was 5277.5
$ws.Cells.Item(67, 12).Value = 4518.3335  # L67 was 5277.5
$ws.Cells.Item(67, 14).Value = -6234.3335  # N67 was -6993.5
$ws.Cells.Item(103, 8).Value = 9113  # H103 was 10050.223
$ws.Cells.Item(103, 9).Value = 681.1111  # I103 was 681.5
$ws.Cells.Item(103, 11).Value = 2043.3333  # K103 was 2044.5
$ws.Cells.Item(103, 13).Value = -1457.3333  # M103 was -1458.5
$ws.Cells.Item(112, 8).Value = 459112.06  # H112 was 494335.03
$ws.Cells.Item(112, 10).Value = 507353.6  # J112 was 550737
$ws.Cells.Item(112, 12).Value = 1522060.8  # L112 was 1652211
$ws.Cells.Item(112, 14).Value = -1524276.8  # N112 was -1654427
$ws.Cells.Item(118, 8).Value = 657.5  # H118 was 747.2727
$ws.Cells.Item(118, 9).Value = 381.42856  # I118 was 418.2
$ws.Cells.Item(118, 10).Value = 1044  # J118 was 1021.5
$ws.Cells.Item(118, 11).Value = 1144.28568  # K118 was 1254.6
$ws.Cells.Item(118, 12).Value = 3132  # L118 was 3064.5
$ws.Cells.Item(118, 13).Value = 512.71432  # M118 was 402.4000000000001
$ws.Cells.Item(118, 14).Value = -6446  # N118 was -6378.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(43, 8).Value = 27441.5  # H43 was 27444
$ws.Cells.Item(43, 10).Value = 27441.5  # J43 was 27444
$ws.Cells.Item(43, 12).Value = 27441.5  # L43 was 27444
$ws.Cells.Item(43, 14).Value = -28067.5  # N43 was -28070
$ws.Cells.Item(63, 8).Value = 9237828  # H63 was 11545784
$ws.Cells.Item(63, 10).Value = 5733.3335  # J63 was 5466.6665
$ws.Cells.Item(63, 12).Value = 5733.3335  # L63 was 5466.6665
$ws.Cells.Item(63, 14).Value = -7105.3335  # N63 was -6838.6665
$ws.Cells.Item(66, 8).Value = 9237828  # H66 was 11545784
$ws.Cells.Item(66, 10).Value = 5733.3335  # J66 was 5466.6665
$ws.Cells.Item(66, 12).Value = 28666.6675  # L66 was 27333.3325
$ws.Cells.Item(66, 14).Value = -35530.6675  # N66 was -34197.3325
$ws.Cells.Item(74, 8).Value = 6020.278  # H74 was 5488.35
$ws.Cells.Item(74, 9).Value = 6488.846  # I74 was 5717.1333
$ws.Cells.Item(74, 11).Value = 6488.846  # K74 was 5717.1333
$ws.Cells.Item(74, 13).Value = -5614.846  # M74 was -4843.1333
$ws.Cells.Item(77, 8).Value = 6020.278  # H77 was 5488.35
$ws.Cells.Item(77, 9).Value = 6488.846  # I77 was 5717.1333
$ws.Cells.Item(77, 11).Value = 32444.23  # K77 was 28585.6665
$ws.Cells.Item(77, 13).Value = -28076.23  # M77 was -24217.6665
$ws.Cells.Item(110, 8).Value = 1261.3334  # H110 was 1180.1904
$ws.Cells.Item(110, 9).Value = 1261.3334  # I110 was 1215.3158
$ws.Cells.Item(110, 10).Value = 0  # J110 was 846.5
$ws.Cells.Item(110, 11).Value = 1261.3334  # K110 was 1215.3158
$ws.Cells.Item(110, 12).Value = 0  # L110 was 846.5
$ws.Cells.Item(110, 13).Value = 783.6666  # M110 was 829.6841999999999
$ws.Cells.Item(110, 14).ClearContents()  # N110 was -4936.5, removed
$ws.Cells.Item(139, 8).Value = 43229.69  # H139 was 43262.77
$ws.Cells.Item(139, 10).Value = 43229.69  # J139 was 43262.77
$ws.Cells.Item(139, 12).Value = 43229.69  # L139 was 43262.77
$ws.Cells.Item(139, 14).Value = -53509.69  # N139 was -53542.77

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 7674.7393  # H20 was 6664.08
$ws.Cells.Item(20, 9).Value = 1499.4546  # I20 was 1105.8462
$ws.Cells.Item(20, 10).Value = 13335.417  # J20 was 12685.5
$ws.Cells.Item(20, 11).Value = 1499.4546  # K20 was 1105.8462
$ws.Cells.Item(20, 12).Value = 13335.417  # L20 was 12685.5
$ws.Cells.Item(20, 13).Value = -1252.4546  # M20 was -858.8462
$ws.Cells.Item(20, 14).Value = -13829.417  # N20 was -13179.5
$ws.Cells.Item(138, 8).Value = 41454.547  # H138 was 41354.074
$ws.Cells.Item(138, 10).Value = 41454.547  # J138 was 41354.074
$ws.Cells.Item(138, 12).Value = 41454.547  # L138 was 41354.074
$ws.Cells.Item(138, 14).Value = -51734.547  # N138 was -51634.074

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(9, 8).Value = 20950  # H9 was 18758
$ws.Cells.Item(9, 10).Value = 20950  # J9 was 18758
$ws.Cells.Item(9, 12).Value = 20950  # L9 was 18758
$ws.Cells.Item(9, 14).Value = -21286  # N9 was -19094
$ws.Cells.Item(16, 8).Value = 10102375  # H16 was 6537053.5
$ws.Cells.Item(16, 9).Value = 13890015  # I16 was 10101901
$ws.Cells.Item(16, 10).Value = 2000  # J16 was 1500
$ws.Cells.Item(16, 11).Value = 13890015  # K16 was 10101901
$ws.Cells.Item(16, 12).Value = 2000  # L16 was 1500
$ws.Cells.Item(16, 13).Value = -13889728  # M16 was -10101614
$ws.Cells.Item(16, 14).Value = -2574  # N16 was -2074
$ws.Cells.Item(41, 8).Value = 31476  # H41 was 43813
$ws.Cells.Item(41, 9).Value = 800  # I41 was 0
$ws.Cells.Item(41, 10).Value = 35858.285  # J41 was 43813
$ws.Cells.Item(41, 11).Value = 800  # K41 was 0
$ws.Cells.Item(41, 12).Value = 35858.285  # L41 was 43813
$ws.Cells.Item(41, 13).Value = -372  # M41 was None
$ws.Cells.Item(41, 14).Value = -36714.285  # N41 was -44669
$ws.Cells.Item(59, 8).Value = 30099.8  # H59 was 32322.572
$ws.Cells.Item(59, 10).Value = 30099.8  # J59 was 32322.572
$ws.Cells.Item(59, 12).Value = 30099.8  # L59 was 32322.572
$ws.Cells.Item(59, 14).Value = -32389.8  # N59 was -34612.572
$ws.Cells.Item(60, 8).Value = 23093.176  # H60 was 25323.455
$ws.Cells.Item(60, 10).Value = 23093.176  # J60 was 25323.455
$ws.Cells.Item(60, 12).Value = 23093.176  # L60 was 25323.455
$ws.Cells.Item(60, 14).Value = -24115.176  # N60 was -26345.455
$ws.Cells.Item(68, 8).Value = 47650.9  # H68 was 45591.727
$ws.Cells.Item(68, 10).Value = 47650.9  # J68 was 45591.727
$ws.Cells.Item(68, 12).Value = 47650.9  # L68 was 45591.727
$ws.Cells.Item(68, 14).Value = -49148.9  # N68 was -47089.727
$ws.Cells.Item(71, 8).Value = 47650.9  # H71 was 45591.727
$ws.Cells.Item(71, 10).Value = 47650.9  # J71 was 45591.727
$ws.Cells.Item(71, 12).Value = 142952.7  # L71 was 136775.181
$ws.Cells.Item(71, 14).Value = -150440.7  # N71 was -144263.181
$ws.Cells.Item(74, 8).Value = 31740  # H74 was 31276.572
$ws.Cells.Item(74, 10).Value = 31740  # J74 was 31276.572
$ws.Cells.Item(74, 12).Value = 31740  # L74 was 31276.572
$ws.Cells.Item(74, 14).Value = -33488  # N74 was -33024.572
$ws.Cells.Item(77, 8).Value = 31740  # H77 was 31276.572
$ws.Cells.Item(77, 10).Value = 31740  # J77 was 31276.572
$ws.Cells.Item(77, 12).Value = 95220  # L77 was 93829.716
$ws.Cells.Item(77, 14).Value = -103956  # N77 was -102565.716
$ws.Cells.Item(105, 8).Value = 2287.4167  # H105 was 3003.1667
$ws.Cells.Item(105, 9).Value = 1944.9  # I105 was 2504.75
$ws.Cells.Item(105, 11).Value = 1944.9  # K105 was 2504.75
$ws.Cells.Item(105, 13).Value = -197.9000000000001  # M105 was -757.75
$ws.Cells.Item(106, 8).Value = 29438.715  # H106 was 31945.166
$ws.Cells.Item(106, 10).Value = 29438.715  # J106 was 31945.166
$ws.Cells.Item(106, 12).Value = 29438.715  # L106 was 31945.166
$ws.Cells.Item(106, 14).Value = -31962.715  # N106 was -34469.166
$ws.Cells.Item(113, 8).Value = 10102375  # H113 was 6537053.5
$ws.Cells.Item(113, 9).Value = 13890015  # I113 was 10101901
$ws.Cells.Item(113, 10).Value = 2000  # J113 was 1500
$ws.Cells.Item(113, 11).Value = 13890015  # K113 was 10101901
$ws.Cells.Item(113, 12).Value = 2000  # L113 was 1500
$ws.Cells.Item(113, 13).Value = -13887845  # M113 was -10099731
$ws.Cells.Item(113, 14).Value = -6340  # N113 was -5840
$ws.Cells.Item(138, 8).Value = 42864.445  # H138 was 44397.145
$ws.Cells.Item(138, 10).Value = 42864.445  # J138 was 44397.145
$ws.Cells.Item(138, 12).Value = 42864.445  # L138 was 44397.145
$ws.Cells.Item(138, 14).Value = -53144.445  # N138 was -54677.145
$ws.Cells.Item(140, 8).Value = 79950.64  # H140 was 79618.5
$ws.Cells.Item(140, 10).Value = 79950.64  # J140 was 79618.5
$ws.Cells.Item(140, 12).Value = 79950.64  # L140 was 79618.5
$ws.Cells.Item(140, 14).Value = -90310.64  # N140 was -89978.5
$ws.Cells.Item(141, 8).Value = 31066.666  # H141 was 30119.2
$ws.Cells.Item(141, 10).Value = 31066.666  # J141 was 30119.2
$ws.Cells.Item(141, 12).Value = 31066.666  # L141 was 30119.2
$ws.Cells.Item(141, 14).Value = -41426.666  # N141 was -40479.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 660  # H113 was 553.43475
$ws.Cells.Item(113, 9).Value = 616.2  # I113 was 526.5
$ws.Cells.Item(113, 11).Value = 1848.6  # K113 was 1579.5
$ws.Cells.Item(113, 13).Value = 321.3999999999999  # M113 was 590.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(140, 8).Value = 37775.383  # H140 was 38518.57
$ws.Cells.Item(140, 10).Value = 37775.383  # J140 was 38518.57
$ws.Cells.Item(140, 12).Value = 37775.383  # L140 was 38518.57
$ws.Cells.Item(140, 14).Value = -48135.383  # N140 was -48878.57

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(139, 8).Value = 44995  # H139 was 43096.25
$ws.Cells.Item(139, 10).Value = 44995  # J139 was 43096.25
$ws.Cells.Item(139, 12).Value = 44995  # L139 was 43096.25
$ws.Cells.Item(139, 14).Value = -55275  # N139 was -53376.25
$ws.Cells.Item(140, 8).Value = 71684.38  # H140 was 73658.086
$ws.Cells.Item(140, 10).Value = 71684.38  # J140 was 73658.086
$ws.Cells.Item(140, 12).Value = 71684.38  # L140 was 73658.086
$ws.Cells.Item(140, 14).Value = -82044.38  # N140 was -84018.086
$ws.Cells.Item(141, 8).Value = 39208.184  # H141 was 40053.637
$ws.Cells.Item(141, 10).Value = 39208.184  # J141 was 40053.637
$ws.Cells.Item(141, 12).Value = 39208.184  # L141 was 40053.637
$ws.Cells.Item(141, 14).Value = -49568.184  # N141 was -50413.637

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 11910154  # H132 was 7095569.5
$ws.Cells.Item(132, 9).Value = 5643.4546  # I132 was 3253.225
$ws.Cells.Item(132, 10).Value = 55560024  # J132 was 47623092
$ws.Cells.Item(132, 11).Value = 16930.3638  # K132 was 9759.674999999999
$ws.Cells.Item(132, 12).Value = 166680072  # L132 was 142869276
$ws.Cells.Item(132, 13).Value = -14400.3638  # M132 was -7229.674999999999
$ws.Cells.Item(132, 14).Value = -166685132  # N132 was -142874336
$ws.Cells.Item(138, 8).Value = 42599.6  # H138 was 40649.668
$ws.Cells.Item(138, 10).Value = 42599.6  # J138 was 40649.668
$ws.Cells.Item(138, 12).Value = 42599.6  # L138 was 40649.668
$ws.Cells.Item(138, 14).Value = -52879.6  # N138 was -50929.668
$ws.Cells.Item(139, 8).Value = 37159.316  # H139 was 37049.617
$ws.Cells.Item(139, 10).Value = 37159.316  # J139 was 37049.617
$ws.Cells.Item(139, 12).Value = 37159.316  # L139 was 37049.617
$ws.Cells.Item(139, 14).Value = -47439.316  # N139 was -47329.617
$ws.Cells.Item(140, 8).Value = 33685.43  # H140 was 33022
$ws.Cells.Item(140, 10).Value = 33685.43  # J140 was 33022
$ws.Cells.Item(140, 12).Value = 33685.43  # L140 was 33022
$ws.Cells.Item(140, 14).Value = -44045.43  # N140 was -43382
$ws.Cells.Item(141, 8).Value = 42392.918  # H141 was 42070.383
$ws.Cells.Item(141, 10).Value = 42392.918  # J141 was 42070.383
$ws.Cells.Item(141, 12).Value = 42392.918  # L141 was 42070.383
$ws.Cells.Item(141, 14).Value = -52752.918  # N141 was -52430.383
